$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 6546
$ws1.Range("F6").Value = 5317
$ws1.Range("F11").Value = 228
$ws1.Range("F12").Value = 39

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 6546
$ws4.Range("F6").Value = 5317
$ws4.Range("F11").Value = 228
$ws4.Range("F14").Value = 39
